$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Late")
$ws2 = $wb.Worksheets.Item("Early Leave")

# New attendance entries logged on the "Late" sheet
$ws1.Range("A4").Value = "13/10/2018 22:41"
$ws1.Range("B4").Value = "T0121212Y"
$ws1.Range("A5").Value = "13/10/2018 22:41"
$ws1.Range("B5").Value = "T0121212Y"
$ws1.Range("A6").ClearFormats()
$ws1.Range("A6").Value = "13/10/2018 22:45"
$ws1.Range("B6").Value = "TESTING1"

# New attendance entries logged on the "Early Leave" sheet
$ws2.Range("A19").Value = "13/10/2018 22:39"
$ws2.Range("B19").Value = "T0121212Y"
$ws2.Range("A20").Value = "31/10/2018 20:23"
$ws2.Range("B20").Value = "asd"

# The "Late" sheet becomes the active/selected sheet (was "Early Leave")
$ws1.Activate()
$ws1.Range("E4").Select()
